$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 440, shifting existing rows 440:489 down to 441:490
$ws.Rows.Item(440).Insert()

# Populate the newly inserted row 440 with the new record
$ws.Range("A440").Value = 5
$ws.Range("B440").Value = "Macroferia Regional de Talca"
$ws.Range("C440").Value = "Maule"
$ws.Range("D440").Value = 44946
$ws.Range("E440").Value = 7
$ws.Range("F440").Value = 100112032
$ws.Range("G440").Value = "Zapallo italiano"
$ws.Range("H440").Value = "Sin especificar"
$ws.Range("I440").Value = "Primera"
$ws.Range("J440").Value = 1200
$ws.Range("K440").Value = 6000
$ws.Range("L440").Value = 6000
$ws.Range("M440").Value = 6000
$ws.Range("N440").Value = "`$/caja 50 unidades"
$ws.Range("O440").Value = "Región del Maule"
$ws.Range("P440").Value = 120
$ws.Range("Q440").Value = 50
$ws.Range("R440").Value = "Hortaliza"
